$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column V labels / values added, plus the U6 "s" -> "m/s" fix.
# Cells are written in the order that makes new shared strings land at
# the same indices as the target workbook (new strings are appended to
# the shared-string table in first-use order).

# Row 3: "?check"  (new shared string)
$ws.Range("V3").Value = "?check"

# Row 2: "protocol"  (new shared string)
$ws.Range("V2").Value = "protocol"

# Row 6: "average"  (new shared string) + fix U6 cm/s/s typo -> m/s
$ws.Range("V6").Value = "average"
$ws.Range("U6").Value = "m/s"

# Row 4: "max"  (new shared string)
$ws.Range("V4").Value = "max"

# Row 7: "?"  (new shared string)
$ws.Range("V7").Value = "?"

# Row 8: "average of max??"  (new shared string)
$ws.Range("V8").Value = "average of max??"

# Row 10: "averageof max from each hand?"  (new shared string)
$ws.Range("V10").Value = "averageof max from each hand?"

# Remaining V cells reuse already-registered shared strings.
$ws.Range("V5").Value = "max"
$ws.Range("V9").Value = "average"
$ws.Range("V11").Value = "max"

# Move the active selection to V12 (matches the saved view state).
[void]$ws.Range("V12").Select()
